# Update NATMI TPM-derived statistics for the Fgf22-Fgfrl1 LR-pair sheet.
#
# The source data (NATMI "lrc2p" output) was regenerated with a new TPM
# normalisation, which changes the ligand/receptor expression-derived
# specificity scores and downstream edge weights in columns E:T for every
# sending/target cluster combination on the sheet. Sending/ligand/receptor/
# target-cluster labels (columns A:D) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.4072979538991744
$ws.Range("J2").Value = 0.4072979538991744
$ws.Range("M2").Value = 1.921622333333333
$ws.Range("N2").Value = 5.764867
$ws.Range("O2").Value = 0.1392241219313625
$ws.Range("P2").Value = 0.1392241219313625
$ws.Range("Q2").Value = 0.2103728076455555
$ws.Range("R2").Value = 1.89335526881
$ws.Range("S2").Value = 0.05670569999605313
$ws.Range("T2").Value = 0.05670569999605313

# Row 3
$ws.Range("I3").Value = 0.4072979538991744
$ws.Range("J3").Value = 0.4072979538991744
$ws.Range("O3").Value = 0.7511588049189343
$ws.Range("P3").Value = 0.7511588049189343
$ws.Range("S3").Value = 0.3059454442968311
$ws.Range("T3").Value = 0.305945444296831

# Row 4
$ws.Range("I4").Value = 0.4072979538991744
$ws.Range("J4").Value = 0.4072979538991744
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1441973333333333
$ws.Range("N4").Value = 0.432592
$ws.Range("O4").Value = 0.01044729069283506
$ws.Range("P4").Value = 0.01044729069283506
$ws.Range("Q4").Value = 0.01578624339555555
$ws.Range("R4").Value = 0.14207619056
$ws.Range("S4").Value = 0.004255160122981608
$ws.Range("T4").Value = 0.004255160122981608

# Row 5
$ws.Range("I5").Value = 0.4072979538991744
$ws.Range("J5").Value = 0.4072979538991744
$ws.Range("M5").Value = 1.182384
$ws.Range("N5").Value = 3.547152
$ws.Range("O5").Value = 0.08566531067535062
$ws.Range("P5").Value = 0.08566531067535062
$ws.Range("Q5").Value = 0.12944345904
$ws.Range("R5").Value = 1.16499113136
$ws.Range("S5").Value = 0.03489130575820741
$ws.Range("T5").Value = 0.03489130575820741

# Row 6
$ws.Range("I6").Value = 0.4072979538991744
$ws.Range("J6").Value = 0.4072979538991744
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1863936666666667
$ws.Range("N6").Value = 0.559181
$ws.Range("O6").Value = 0.01350447178151746
$ws.Range("P6").Value = 0.01350447178151746
$ws.Range("Q6").Value = 0.02040575731444445
$ws.Range("R6").Value = 0.18365181583
$ws.Range("S6").Value = 0.005500343725101202
$ws.Range("T6").Value = 0.005500343725101202

# Row 7
$ws.Range("G7").Value = 0.09509666666666666
$ws.Range("H7").Value = 0.28529
$ws.Range("I7").Value = 0.3537984753764744
$ws.Range("J7").Value = 0.3537984753764744
$ws.Range("M7").Value = 1.921622333333333
$ws.Range("N7").Value = 5.764867
$ws.Range("O7").Value = 0.1392241219313625
$ws.Range("P7").Value = 0.1392241219313625
$ws.Range("Q7").Value = 0.1827398784922222
$ws.Range("R7").Value = 1.64465890643
$ws.Range("S7").Value = 0.04925728207494443
$ws.Range("T7").Value = 0.04925728207494443

# Row 8
$ws.Range("G8").Value = 0.09509666666666666
$ws.Range("H8").Value = 0.28529
$ws.Range("I8").Value = 0.3537984753764744
$ws.Range("J8").Value = 0.3537984753764744
$ws.Range("O8").Value = 0.7511588049189343
$ws.Range("P8").Value = 0.7511588049189343
$ws.Range("Q8").Value = 0.9859402726700001
$ws.Range("R8").Value = 8.873462454029999
$ws.Range("S8").Value = 0.2657588399459335
$ws.Range("T8").Value = 0.2657588399459335

# Row 9
$ws.Range("G9").Value = 0.09509666666666666
$ws.Range("H9").Value = 0.28529
$ws.Range("I9").Value = 0.3537984753764744
$ws.Range("J9").Value = 0.3537984753764744
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1441973333333333
$ws.Range("N9").Value = 0.432592
$ws.Range("O9").Value = 0.01044729069283506
$ws.Range("P9").Value = 0.01044729069283506
$ws.Range("Q9").Value = 0.01371268574222222
$ws.Range("R9").Value = 0.12341417168
$ws.Range("S9").Value = 0.003696235518939875
$ws.Range("T9").Value = 0.003696235518939875

# Row 10
$ws.Range("G10").Value = 0.09509666666666666
$ws.Range("H10").Value = 0.28529
$ws.Range("I10").Value = 0.3537984753764744
$ws.Range("J10").Value = 0.3537984753764744
$ws.Range("M10").Value = 1.182384
$ws.Range("N10").Value = 3.547152
$ws.Range("O10").Value = 0.08566531067535062
$ws.Range("P10").Value = 0.08566531067535062
$ws.Range("Q10").Value = 0.11244077712
$ws.Range("R10").Value = 1.01196699408
$ws.Range("S10").Value = 0.03030825630959106
$ws.Range("T10").Value = 0.03030825630959106

# Row 11
$ws.Range("G11").Value = 0.09509666666666666
$ws.Range("H11").Value = 0.28529
$ws.Range("I11").Value = 0.3537984753764744
$ws.Range("J11").Value = 0.3537984753764744
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1863936666666667
$ws.Range("N11").Value = 0.559181
$ws.Range("O11").Value = 0.01350447178151746
$ws.Range("P11").Value = 0.01350447178151746
$ws.Range("Q11").Value = 0.01772541638777778
$ws.Range("R11").Value = 0.15952874749
$ws.Range("S11").Value = 0.004777861527065498
$ws.Range("T11").Value = 0.004777861527065499

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.06421433333333333
$ws.Range("H12").Value = 0.192643
$ws.Range("I12").Value = 0.2389035707243512
$ws.Range("J12").Value = 0.2389035707243512
$ws.Range("M12").Value = 1.921622333333333
$ws.Range("N12").Value = 5.764867
$ws.Range("O12").Value = 0.1392241219313625
$ws.Range("P12").Value = 0.1392241219313625
$ws.Range("Q12").Value = 0.1233956970534444
$ws.Range("R12").Value = 1.110561273481
$ws.Range("S12").Value = 0.03326113986036496
$ws.Range("T12").Value = 0.03326113986036496

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.06421433333333333
$ws.Range("H13").Value = 0.192643
$ws.Range("I13").Value = 0.2389035707243512
$ws.Range("J13").Value = 0.2389035707243512
$ws.Range("O13").Value = 0.7511588049189343
$ws.Range("P13").Value = 0.7511588049189343
$ws.Range("Q13").Value = 0.665759374489
$ws.Range("R13").Value = 5.991834370401
$ws.Range("S13").Value = 0.1794545206761697
$ws.Range("T13").Value = 0.1794545206761697

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.06421433333333333
$ws.Range("H14").Value = 0.192643
$ws.Range("I14").Value = 0.2389035707243512
$ws.Range("J14").Value = 0.2389035707243512
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1441973333333333
$ws.Range("N14").Value = 0.432592
$ws.Range("O14").Value = 0.01044729069283506
$ws.Range("P14").Value = 0.01044729069283506
$ws.Range("Q14").Value = 0.009259535628444443
$ws.Range("R14").Value = 0.083335820656
$ws.Range("S14").Value = 0.002495895050913577
$ws.Range("T14").Value = 0.002495895050913577

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.06421433333333333
$ws.Range("H15").Value = 0.192643
$ws.Range("I15").Value = 0.2389035707243512
$ws.Range("J15").Value = 0.2389035707243512
$ws.Range("M15").Value = 1.182384
$ws.Range("N15").Value = 3.547152
$ws.Range("O15").Value = 0.08566531067535062
$ws.Range("P15").Value = 0.08566531067535062
$ws.Range("Q15").Value = 0.075926000304
$ws.Range("R15").Value = 0.6833340027360001
$ws.Range("S15").Value = 0.02046574860755214
$ws.Range("T15").Value = 0.02046574860755214

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.06421433333333333
$ws.Range("H16").Value = 0.192643
$ws.Range("I16").Value = 0.2389035707243512
$ws.Range("J16").Value = 0.2389035707243512
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1863936666666667
$ws.Range("N16").Value = 0.559181
$ws.Range("O16").Value = 0.01350447178151746
$ws.Range("P16").Value = 0.01350447178151746
$ws.Range("Q16").Value = 0.01196914504255556
$ws.Range("R16").Value = 0.107722305383
$ws.Range("S16").Value = 0.003226266529350762
$ws.Range("T16").Value = 0.003226266529350763
